$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 100, shifting all
# subsequent rows (old 100..174) down to 102..176.
$ws.Rows.Item(100).Insert()
$ws.Rows.Item(100).Insert()

# New row 100 data
$ws.Range("A100").Value = 5
$ws.Range("B100").Value = "Macroferia Regional de Talca"
$ws.Range("C100").Value = "Maule"
$ws.Range("D100").Value = 45216
$ws.Range("E100").Value = 7
$ws.Range("F100").Value = "Fruta"
$ws.Range("G100").Value = 100107
$ws.Range("H100").Value = "Otros"
$ws.Range("I100").Value = 100107002
$ws.Range("J100").Value = "Chirimoya"
$ws.Range("K100").Value = "Cultivar IV Región"
$ws.Range("L100").Value = "Primera"
$ws.Range("M100").Value = 120
$ws.Range("N100").Value = 22000
$ws.Range("O100").Value = 22000
$ws.Range("P100").Value = 22000
$ws.Range("Q100").Value = "$/bandeja 10 kilos"
$ws.Range("R100").Value = "Provincia de Limarí"
$ws.Range("S100").Value = 2200
$ws.Range("T100").Value = 10

# New row 101 data
$ws.Range("A101").Value = 5
$ws.Range("B101").Value = "Macroferia Regional de Talca"
$ws.Range("C101").Value = "Maule"
$ws.Range("D101").Value = 45216
$ws.Range("E101").Value = 7
$ws.Range("F101").Value = "Fruta"
$ws.Range("G101").Value = 100107
$ws.Range("H101").Value = "Otros"
$ws.Range("I101").Value = 100107002
$ws.Range("J101").Value = "Chirimoya"
$ws.Range("K101").Value = "Cultivar IV Región"
$ws.Range("L101").Value = "Segunda"
$ws.Range("M101").Value = 180
$ws.Range("N101").Value = 20000
$ws.Range("O101").Value = 20000
$ws.Range("P101").Value = 20000
$ws.Range("Q101").Value = "$/bandeja 10 kilos"
$ws.Range("R101").Value = "Provincia de Limarí"
$ws.Range("S101").Value = 2000
$ws.Range("T101").Value = 10
